$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to write (crypto Price/Volume(1h) columns).
$updates = [ordered]@{
    'D2' = '62.386.16'
    'D3' = '2.422.42'
    'E3' = '  +0.70%  '
    'E4' = '  -0.02%  '
    'D5' = '564.12'
    'E5' = '  +0.49%  '
    'D6' = '144.54'
    'E7' = '  +0.00%  '
    'E8' = '  +0.24%  '
    'E9' = '  +1.31%  '
    'E11' = '  +1.00%  '
    'D12' = '0.353'
    'E12' = '  +1.67%  '
    'D13' = '26.70'
    'E13' = '  +4.99%  '
    'D14' = '0.0000178'
    'D15' = '2.860.62'
    'E15' = '  +0.69%  '
    'D16' = '62.269.48'
    'E16' = '  +0.76%  '
    'D17' = '2.449.69'
    'E17' = '  +1.85%  '
    'D18' = '11.17'
    'E18' = '  -0.20%  '
    'D19' = '6.94'
    'E19' = '  +2.02%  '
    'D20' = '322.69'
    'E20' = '  +0.66%  '
    'E21' = '  +1.03%  '
    'D22' = '1.00'
    'E22' = '  -0.01%  '
    'D23' = '67.03'
    'E23' = '  +2.17%  '
    'D24' = '1.81'
    'E24' = '  +4.64%  '
    'D25' = '586.57'
    'E25' = '  +4.14%  '
    'D26' = '8.52'
    'E26' = '  -1.56%  '
    'E27' = '  +7.49%  '
    'D28' = '2.540.15'
    'E28' = '  +0.99%  '
    'D29' = '1.00'
    'E29' = '  +0.04%  '
    'D30' = '8.39'
    'E30' = '  +2.99%  '
    'D31' = '1.43'
    'E31' = '  +3.81%  '
    'D32' = '0.145'
    'E32' = '  -1.36%  '
    'E33' = '  +0.25%  '
    'E34' = '  -0.39%  '
    'E35' = '  +0.07%  '
    'D36' = '4.81'
    'E36' = '  +1.63%  '
    'D37' = '0.380'
    'E37' = '  +0.41%  '
    'D38' = '18.67'
    'E38' = '  +1.13%  '
    'D39' = '5.32'
    'E39' = '  -1.74%  '
    'D40' = '147.86'
    'E40' = '  -2.87%  '
    'E41' = '  +1.64%  '
    'D42' = '0.999'
    'E42' = '  +0.02%  '
    'D43' = '2.43'
    'E43' = '  +9.26%  '
    'D44' = '147.62'
    'D45' = '3.65'
    'E45' = '  +1.96%  '
    'D46' = '0.0533'
    'E46' = '  +1.03%  '
    'D47' = '20.41'
    'E47' = '  +3.38%  '
    'D48' = '0.598'
    'E48' = '  +1.74%  '
    'D49' = '0.0230'
    'E49' = '  +2.76%  '
    'D50' = '0.0917'
    'E50' = '  +0.08%  '
    'E51' = '  +4.26%  '
}

# Price-column cells whose new text looks like a plain number (e.g. "1.00",
# "0.0230") must be forced to Text format first, otherwise Excel COM auto-
# converts the assigned string into a numeric value and the literal
# formatting (trailing zeros, etc.) is lost.
$forceText = @(
    'D5', 'D6', 'D12', 'D13', 'D14', 'D18', 'D19', 'D20',
    'D22', 'D23', 'D24', 'D25', 'D26', 'D29', 'D30', 'D31',
    'D32', 'D36', 'D37', 'D38', 'D39', 'D40', 'D42', 'D43',
    'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50'
)

foreach ($addr in $forceText) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells ($($forceText.Count) forced to text)"
